# Add two new API rows (rescue endpoints) to the bottom of the API table,
# matching the new unit test cases added to the spec sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: GET /api/rescue/requests -> list of rescues
$ws.Range("A21").Value = "/api/rescue/requests"
$ws.Range("B21").Value = "get"
$ws.Range("C21").Value = "list of rescues"

# Row 22: POST /api/rescue/request -> register for rescue request
# (fill description/method before the endpoint column so the shared-string
# table ends up in the same order as a natural left-to-right-ish entry)
$ws.Range("C22").Value = "register for rescue request"
$ws.Range("B22").Value = "post"
$ws.Range("A22").Value = "/api/rescue/request"

# Column A entries use the 12pt "endpoint" font/style used throughout the
# rest of the table, and every data row uses the 15.75pt row height.
$ws.Range("A21:A22").Font.Size = 12
$ws.Range("A21:C22").RowHeight = 15.75

# Scroll the view down and leave the new second row selected, mirroring
# where the author ended up after typing in the new rows.
$ws.Application.ActiveWindow.ScrollRow = 15
$ws.Range("A22").Select()
